# Bugfix in zve(): cap the "Zuschuss/Verguenstigung" term (12*(P+N+0.96*O))
# at 1900 using MIN(), for every row of the Z column (rows 2-25).
# Row 13 keeps its own special-cased formula (hard-coded 0.0995*5500 instead
# of M13), so it is written separately to preserve its distinct formula text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 has its own (non-shared) formula.
$ws.Range("Z2").Formula = "=((0.6+(0.02*(T2-2005)))*(2*12*M2))-(12*M2)+MIN(12*(P2+N2+0.96*O2), 1900)"

# Rows 3-12 share one formula pattern (relative refs adjust per row).
$ws.Range("Z3:Z12").Formula = "=((0.6+(0.02*(T3-2005)))*(2*12*M3))-(12*M3)+MIN(12*(P3+N3+0.96*O3), 1900)"

# Row 13 is special: it hard-codes the pension contribution (0.0995*5500)
# instead of referencing M13.
$ws.Range("Z13").Formula = "=((0.6+(0.02*(T13-2005)))*(12*2*0.0995*5500))-(12*0.0995*5500)+MIN(12*(P13+N13+0.96*O13),1900)"

# Rows 14-25 continue the same shared pattern as rows 3-12.
$ws.Range("Z14:Z25").Formula = "=((0.6+(0.02*(T14-2005)))*(2*12*M14))-(12*M14)+MIN(12*(P14+N14+0.96*O14), 1900)"

# Reflect the workbook's last active selection at save time.
$ws.Range("Z13").Select()
